# Vragen Tech Assistant - content edits
$d = $word.ActiveDocument

# 1) Replace "Welke versie van Python" question with the new datastructuur question
$d.Content.Find.Execute("Welke versie van Python", $false, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Wat is de beste datastructuur om de stations met coördinaten en de trajecten met reistijden in te laden?",
                         2)

# 2) Delete the "Welke ontwikkelomgeving" paragraph entirely (including its paragraph mark)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Welke ontwikkelomgeving*") {
        $target = $p
    }
}
if ($target) { $target.Range.Delete() }

# 3) Delete the "Hoe werk je tegelijkertijd..." paragraph entirely
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Hoe werk je tegelijkertijd*") {
        $target2 = $p
    }
}
if ($target2) { $target2.Range.Delete() }

# 4) array -> list in "Tweedimensionale array met coördinaten..."
$d.Content.Find.Execute("Tweedimensionale array met coördinaten", $false, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Tweedimensionale list met coördinaten",
                         2)

# 5) Array -> List in "Array op basis van afstand..."
$d.Content.Find.Execute("Array op basis van afstand", $false, $false, $false, $false, $false,
                         $true, 1, $false,
                         "List op basis van afstand",
                         2)
